# Applies the "objective_binge_single_cat" config change / bug fix:
#   - H5/I5, H9/I9, H10/I10 used to hold a shared formula
#     (=IFERROR(1/F,2) / =IFERROR(1/G,2)); they are replaced with plain,
#     hard-coded numeric values.
#   - The active sheet's view/selection moves from I18 to O10 (scrolled so
#     row 3 is the top visible row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bug fix: bake the (previously shared/formula-driven) ratios into
#     plain literal values ------------------------------------------------
$ws.Range("H5").Value  = 0.25
$ws.Range("I5").Value  = 0.25

$ws.Range("H9").Value  = 0.25
$ws.Range("I9").Value  = 0.5

$ws.Range("H10").Value = 0.1666
$ws.Range("I10").Value = 0.25

# --- View/selection change: scroll so row 3 is at the top and select O10 -
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow    = 3
$win.ScrollColumn = 1
[void]$ws.Range("O10").Select()
